$wb = $excel.ActiveWorkbook

# --- Update "Bedrift" organisasjonsnummer cell (G2): was text "924 88 1 6 82", now numeric 924881682 ---
$bedrift = $wb.Worksheets.Item("Bedrift")
$bedrift.Range("G2").Value = 924881682

# --- Remove the empty, unused "Ark1" worksheet ---
$excel.DisplayAlerts = $false
$ark1 = $wb.Worksheets.Item("Ark1")
$ark1.Delete()
$excel.DisplayAlerts = $true

# --- Re-apply column widths on "Bedrift" that come from the real workbook's layout ---
$bedrift.Columns.Item(2).ColumnWidth = 16.08984375
$bedrift.Columns.Item(3).ColumnWidth = 15.08984375
$bedrift.Columns.Item(4).ColumnWidth = 12.08984375
$bedrift.Columns.Item(5).ColumnWidth = 27.54296875
$bedrift.Columns.Item(7).ColumnWidth = 19.1796875
$bedrift.Columns.Item(8).ColumnWidth = 16

# --- Make "Bedrift" the active sheet/tab, with G2 selected, matching the saved view state ---
$bedrift.Activate()
$bedrift.Range("G2").Select()
